# chore - final commit warnings removed
#
# Refresh the "Created Room ID" values (column L) with the IDs produced by
# the latest live test run against the Booking API. Every other column
# (dataset name, description, expected result, etc.) is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number (column L = "Created Room ID") to its new value.
# These cells are text-typed (General format) numeric strings; we flip the
# NumberFormat to Text ("@") before assigning so Excel keeps the value as a
# string instead of coercing it to a number, then restore "General" so the
# cell style/format matches the original.
$roomIdUpdates = @{
    2 = "435539446"
    9 = "435581115"
    10 = "435582373"
    11 = "435583270"
    12 = "435584363"
    13 = "435585128"
    14 = "435585990"
    15 = "435586801"
    16 = "435588067"
    17 = "435588934"
    18 = "435589824"
    19 = "435590690"
    20 = "435591605"
    21 = "435592556"
    22 = "435593516"
    23 = "435594726"
    24 = "435595955"
    25 = "435596957"
    26 = "435598243"
    27 = "435599209"
    28 = "435600456"
    29 = "435601477"
    30 = "435602784"
    31 = "435603859"
    32 = "435604972"
    33 = "435606106"
    34 = "435607379"
    35 = "435608579"
    36 = "435610016"
    37 = "435611190"
    38 = "435612584"
    39 = "435613878"
    40 = "435615083"
    41 = "435616293"
    42 = "435617520"
    43 = "435618767"
    44 = "435620253"
    45 = "435621519"
    46 = "435622801"
    47 = "435624516"
    48 = "435625915"
    49 = "435627213"
    50 = "435630576"
    51 = "435632100"
    52 = "435633504"
    53 = "435628683"
    54 = "435634950"
    55 = "435636415"
    56 = "435638081"
    57 = "435639590"
    58 = "435641363"
    59 = "435642845"
    60 = "435645310"
    61 = "435674054"
    62 = "435676742"
    63 = "435679813"
    64 = "435683074"
    65 = "435685393"
    66 = "435687952"
    67 = "435690277"
    68 = "435692781"
    69 = "435649070"
    70 = "435659968"
    71 = "435662333"
    72 = "435651776"
    74 = "435656441"
    75 = "435665971"
    76 = "435668629"
    77 = "435671083"
}

foreach ($row in $roomIdUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 12)
    $originalFormat = $cell.NumberFormat
    $cell.NumberFormat = "@"
    $cell.Value = $roomIdUpdates[$row]
    $cell.NumberFormat = $originalFormat
}
